# Auto-generated edit script applying the Brynhildr_Profits market-data refresh
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4290.8335
$ws.Range("J64").Value = 4290.8335
$ws.Range("L64").Value = 4290.8335
$ws.Range("N64").Value = -4786.8335
$ws.Range("H67").Value = 4290.8335
$ws.Range("J67").Value = 4290.8335
$ws.Range("L67").Value = 4290.8335
$ws.Range("N67").Value = -6006.8335
$ws.Range("H116").Value = 17531.785
$ws.Range("J116").Value = 22499.834
$ws.Range("L116").Value = 22499.834
$ws.Range("N116").Value = -29383.834

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 857980.4
$ws.Range("I32").Value = 927317.4
$ws.Range("K32").Value = 927317.4
$ws.Range("M32").Value = -927030.4
$ws.Range("H102").Value = 2012.84
$ws.Range("I102").Value = 1148.4736
$ws.Range("K102").Value = 1148.4736
$ws.Range("M102").Value = 473.5264
$ws.Range("H122").Value = 1823.9615
$ws.Range("I122").Value = 1707.2084
$ws.Range("K122").Value = 5121.6252
$ws.Range("M122").Value = -2671.6252

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43322.76
$ws.Range("I20").Value = 55951.156
$ws.Range("J20").Value = 3332.8333
$ws.Range("K20").Value = 55951.156
$ws.Range("L20").Value = 3332.8333
$ws.Range("M20").Value = -55704.156
$ws.Range("N20").Value = -3826.8333
$ws.Range("H86").Value = 2515.7727
$ws.Range("I86").Value = 2311.1875
$ws.Range("K86").Value = 2311.1875
$ws.Range("M86").Value = -1188.1875
$ws.Range("H89").Value = 2515.7727
$ws.Range("I89").Value = 2311.1875
$ws.Range("K89").Value = 11555.9375
$ws.Range("M89").Value = -5939.9375
$ws.Range("H94").Value = 3056.12
$ws.Range("I94").Value = 1669.7
$ws.Range("J94").Value = 8601.799999999999
$ws.Range("K94").Value = 1669.7
$ws.Range("L94").Value = 8601.799999999999
$ws.Range("M94").Value = -1218.7
$ws.Range("N94").Value = -9503.799999999999
$ws.Range("H105").Value = 1996.027
$ws.Range("I105").Value = 1692.0741
$ws.Range("K105").Value = 1692.0741
$ws.Range("M105").Value = 54.92589999999996

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 8291.076999999999
$ws.Range("J122").Value = 250000
$ws.Range("L122").Value = 750000
$ws.Range("N122").Value = -754900

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 84.166664
$ws.Range("I36").Value = 84.166664
$ws.Range("K36").Value = 252.499992
$ws.Range("M36").Value = -83.49999199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2291.9583
$ws.Range("I80").Value = 2165.2727
$ws.Range("K80").Value = 2165.2727
$ws.Range("M80").Value = -1167.2727
$ws.Range("H83").Value = 2291.9583
$ws.Range("I83").Value = 2165.2727
$ws.Range("K83").Value = 10826.3635
$ws.Range("M83").Value = -5834.363499999999
$ws.Range("H132").Value = 12123
$ws.Range("I132").Value = 12123
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 36369
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -33839
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3508.182
$ws.Range("I46").Value = 1188.6666
$ws.Range("K46").Value = 1188.6666
$ws.Range("M46").Value = -1000.6666
$ws.Range("H68").Value = 7678.2383
$ws.Range("I68").Value = 8836.846
$ws.Range("J68").Value = 5795.5
$ws.Range("K68").Value = 8836.846
$ws.Range("L68").Value = 5795.5
$ws.Range("M68").Value = -8087.846
$ws.Range("N68").Value = -7293.5
$ws.Range("H71").Value = 7678.2383
$ws.Range("I71").Value = 8836.846
$ws.Range("J71").Value = 5795.5
$ws.Range("K71").Value = 44184.23
$ws.Range("L71").Value = 28977.5
$ws.Range("M71").Value = -40440.23
$ws.Range("N71").Value = -36465.5
$ws.Range("H82").Value = 2842.9167
$ws.Range("I82").Value = 2847.5557
$ws.Range("K82").Value = 2847.5557
$ws.Range("M82").Value = -2486.5557
$ws.Range("H85").Value = 2842.9167
$ws.Range("I85").Value = 2847.5557
$ws.Range("K85").Value = 2847.5557
$ws.Range("M85").Value = -1599.5557
$ws.Range("H93").Value = 2102.7
$ws.Range("I93").Value = 1161.125
$ws.Range("K93").Value = 1161.125
$ws.Range("M93").Value = 86.875
$ws.Range("H95").Value = 19397.5
$ws.Range("J95").Value = 19397.5
$ws.Range("L95").Value = 19397.5
$ws.Range("N95").Value = -24889.5
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H100").Value = 2950
$ws.Range("I100").Value = 2950
$ws.Range("K100").Value = 2950
$ws.Range("M100").Value = -2409
$ws.Range("H101").Value = 21000
$ws.Range("J101").Value = 21000
$ws.Range("L101").Value = 21000
$ws.Range("N101").Value = -27490
$ws.Range("H102").Value = 55000
$ws.Range("J102").Value = 55000
$ws.Range("L102").Value = 55000
$ws.Range("N102").Value = -61490
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H106").Value = 27249.5
$ws.Range("J106").Value = 27249.5
$ws.Range("L106").Value = 27249.5
$ws.Range("N106").Value = -29773.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10196.41
$ws.Range("I62").Value = 9989.857
$ws.Range("J62").Value = 10312.08
$ws.Range("K62").Value = 9989.857
$ws.Range("L62").Value = 10312.08
$ws.Range("M62").Value = -9365.857
$ws.Range("N62").Value = -11560.08
$ws.Range("H65").Value = 10196.41
$ws.Range("I65").Value = 9989.857
$ws.Range("J65").Value = 10312.08
$ws.Range("K65").Value = 49949.285
$ws.Range("L65").Value = 51560.4
$ws.Range("M65").Value = -46829.285
$ws.Range("N65").Value = -57800.4
$ws.Range("H81").Value = 103460.5
$ws.Range("I81").Value = 3575.8333
$ws.Range("J81").Value = 253287.5
$ws.Range("K81").Value = 7151.6666
$ws.Range("L81").Value = 506575
$ws.Range("M81").Value = -6090.6666
$ws.Range("N81").Value = -508697
$ws.Range("H84").Value = 103460.5
$ws.Range("I84").Value = 3575.8333
$ws.Range("J84").Value = 253287.5
$ws.Range("K84").Value = 35758.333
$ws.Range("L84").Value = 2532875
$ws.Range("M84").Value = -30454.333
$ws.Range("N84").Value = -2543483
$ws.Range("H122").Value = 32268.053
$ws.Range("I122").Value = 2534.2856
$ws.Range("K122").Value = 7602.8568
$ws.Range("M122").Value = -5152.8568
